# Add a new paragraph after "This an Arduino test file." containing
# "We are testing the file.", inheriting the same paragraph/run formatting
# (Comic Sans MS, bold, orange accent2 color, size 40, centered).

$d = $word.ActiveDocument

$firstPara = $d.Paragraphs.First
$firstPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "We are testing the file."
